# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to leve-profit data cells across all 8 profession sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 6055.9414
$ws.Range("I33").Value = 66.09999999999999
$ws.Range("J33").Value = 14612.857
$ws.Range("K33").Value = 66.09999999999999
$ws.Range("L33").Value = 14612.857
$ws.Range("M33").Value = 162.9
$ws.Range("N33").Value = -15070.857

# Row 51
$ws.Range("H51").Value = 3784
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 3807.6667
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 3807.6667
$ws.Range("M51").Value = -3016
$ws.Range("N51").Value = -4775.6667

# Row 132
$ws.Range("H132").Value = 1753.258
$ws.Range("I132").Value = 1679.7858
$ws.Range("J132").Value = 2439
$ws.Range("K132").Value = 5039.357400000001
$ws.Range("L132").Value = 7317
$ws.Range("M132").Value = -2509.357400000001
$ws.Range("N132").Value = -12377

# Row 140
$ws.Range("H140").Value = 83123
$ws.Range("J140").Value = 83123
$ws.Range("L140").Value = 83123
$ws.Range("N140").Value = -93483

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 311641.66
$ws.Range("I32").Value = 3516.5488
$ws.Range("J32").Value = 2608574.2
$ws.Range("K32").Value = 3516.5488
$ws.Range("L32").Value = 2608574.2
$ws.Range("M32").Value = -3229.5488
$ws.Range("N32").Value = -2609148.2

# Row 37
$ws.Range("H37").Value = 90914230
$ws.Range("I37").Value = 200001600
$ws.Range("K37").Value = 200001600
$ws.Range("M37").Value = -200001327

# Row 61
$ws.Range("H61").Value = 2571.9756
$ws.Range("I61").Value = 2684.5
$ws.Range("J61").Value = 2025.4286
$ws.Range("K61").Value = 2684.5
$ws.Range("L61").Value = 2025.4286
$ws.Range("M61").Value = -2472.5
$ws.Range("N61").Value = -2449.4286

# Row 74
$ws.Range("H74").Value = 2644.4583
$ws.Range("I74").Value = 2541.75
$ws.Range("K74").Value = 2541.75
$ws.Range("M74").Value = -1667.75

# Row 77
$ws.Range("H77").Value = 2644.4583
$ws.Range("I77").Value = 2541.75
$ws.Range("K77").Value = 12708.75
$ws.Range("M77").Value = -8340.75

# Row 136
$ws.Range("H136").Value = 2571.9756
$ws.Range("I136").Value = 2684.5
$ws.Range("J136").Value = 2025.4286
$ws.Range("K136").Value = 8053.5
$ws.Range("L136").Value = 6076.2858
$ws.Range("M136").Value = -5503.5
$ws.Range("N136").Value = -11176.2858

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4402.449
$ws.Range("I134").Value = 1059.4117
$ws.Range("J134").Value = 11980
$ws.Range("K134").Value = 3178.2351
$ws.Range("L134").Value = 35940
$ws.Range("M134").Value = -643.2351000000003
$ws.Range("N134").Value = -41010

# Row 138
$ws.Range("H138").Value = 47450.91
$ws.Range("J138").Value = 47450.91
$ws.Range("L138").Value = 47450.91
$ws.Range("N138").Value = -57730.91

# Row 140
$ws.Range("H140").Value = 61971.43
$ws.Range("J140").Value = 61971.43
$ws.Range("L140").Value = 61971.43
$ws.Range("N140").Value = -72331.42999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 927
$ws.Range("I16").Value = 924.8889
$ws.Range("J16").Value = 933.3333
$ws.Range("K16").Value = 924.8889
$ws.Range("L16").Value = 933.3333
$ws.Range("M16").Value = -637.8889
$ws.Range("N16").Value = -1507.3333

# Row 31
$ws.Range("H31").Value = 4213.316
$ws.Range("I31").Value = 3255.138
$ws.Range("J31").Value = 4804.5317
$ws.Range("K31").Value = 3255.138
$ws.Range("L31").Value = 4804.5317
$ws.Range("M31").Value = -2960.138
$ws.Range("N31").Value = -5394.5317

# Row 34
$ws.Range("H34").Value = 4213.316
$ws.Range("I34").Value = 3255.138
$ws.Range("J34").Value = 4804.5317
$ws.Range("K34").Value = 3255.138
$ws.Range("L34").Value = 4804.5317
$ws.Range("M34").Value = -3053.138
$ws.Range("N34").Value = -5208.5317

# Row 50
$ws.Range("H50").Value = 8575.4
$ws.Range("J50").Value = 8575.4
$ws.Range("L50").Value = 8575.4
$ws.Range("N50").Value = -9825.4

# Row 51
$ws.Range("H51").Value = 8987.875
$ws.Range("J51").Value = 9414.714
$ws.Range("L51").Value = 9414.714
$ws.Range("N51").Value = -10886.714

# Row 60
$ws.Range("H60").Value = 6599.8
$ws.Range("I60").Value = 200
$ws.Range("J60").Value = 8199.75
$ws.Range("K60").Value = 200
$ws.Range("L60").Value = 8199.75
$ws.Range("M60").Value = 311
$ws.Range("N60").Value = -9221.75

# Row 61
$ws.Range("H61").Value = 8987.875
$ws.Range("J61").Value = 9414.714
$ws.Range("L61").Value = 9414.714
$ws.Range("N61").Value = -10110.714

# Row 74
$ws.Range("H74").Value = 14668.3
$ws.Range("J74").Value = 16044.223
$ws.Range("L74").Value = 16044.223
$ws.Range("N74").Value = -17792.223

# Row 77
$ws.Range("H77").Value = 14668.3
$ws.Range("J77").Value = 16044.223
$ws.Range("L77").Value = 48132.669
$ws.Range("N77").Value = -56868.669

# Row 113
$ws.Range("H113").Value = 927
$ws.Range("I113").Value = 924.8889
$ws.Range("J113").Value = 933.3333
$ws.Range("K113").Value = 924.8889
$ws.Range("L113").Value = 933.3333
$ws.Range("M113").Value = 1245.1111
$ws.Range("N113").Value = -5273.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 847.7778
$ws.Range("I109").Value = 440
$ws.Range("J109").Value = 1357.5
$ws.Range("K109").Value = 1320
$ws.Range("L109").Value = 4072.5
$ws.Range("M109").Value = -280
$ws.Range("N109").Value = -6152.5

# Row 131
$ws.Range("H131").Value = 797.1919
$ws.Range("I131").Value = 332
$ws.Range("J131").Value = 821.93616
$ws.Range("K131").Value = 996
$ws.Range("L131").Value = 2465.80848
$ws.Range("M131").Value = 4044
$ws.Range("N131").Value = -12545.80848

# Row 139
$ws.Range("H139").Value = 1815
$ws.Range("I139").Value = 933.6667
$ws.Range("J139").Value = 2635.5518
$ws.Range("K139").Value = 2801.0001
$ws.Range("L139").Value = 7906.655400000001
$ws.Range("M139").Value = 2338.9999
$ws.Range("N139").Value = -18186.6554

$ws = $wb.Worksheets.Item("GSM")
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 666.56366
$ws.Range("I136").Value = 543.1111
$ws.Range("J136").Value = 900.4737
$ws.Range("K136").Value = 1629.3333
$ws.Range("L136").Value = 2701.4211
$ws.Range("M136").Value = 920.6667000000002
$ws.Range("N136").Value = -7801.4211
